$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the description text for the "Project Schedule" CI (row 15)
$ws.Range("D15").Value = "The project schedule is the tool that communicates what work needs to be performed, which resources of the organization will perform the work and the timeframes in which that work needs to be performed."

# Add a new CI row for "work breakdown structure"
$ws.Range("B16").Value = "work breakdown structure"
$ws.Range("C16").Value = "CAR_work breakdown structure"
$ws.Range("D16").Value = "WBS is a key project deliverable that organizes the team's work into manageable sections"
$ws.Range("E16").Value = "Version"

$ws.Rows.Item(16).RowHeight = 47.25

# Extend the Folder merge to cover the new row
$ws.Range("A8:A16").Merge()

$ws.Range("D16").Select()
